# Messreihen.xlsx -- add angle/phase-shift analysis columns (J:O) for rows 2-12,
# header row 13, and a formatted-but-empty K14 cell. Also adjusts column widths
# and the sheet's active selection (matches the "neues matlab analyse programm"
# commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- J column: Kap fuer (20 Grad) ------------------------------------------------
$ws.Range("J2").Formula = "= TAN(20 * 180 / PI()) / (2*PI()*50*G2)"
$ws.Range("J2").NumberFormat = "0.00E+00"
$ws.Range("J3:J12").Formula = "= TAN(20 * 180 / PI()) / (2*PI()*50*G3)"
$ws.Range("J3:J12").NumberFormat = "0.00E+00"

# --- K column: Phi bei 10uF -------------------------------------------------
$ws.Range("K2").Formula = "=ATAN(2*PI()*50*G2*0.00001) * 180 / PI()"
$ws.Range("K2").NumberFormat = "0.00"
$ws.Range("K3:K12").Formula = "=ATAN(2*PI()*50*G3*0.00001) * 180 / PI()"
$ws.Range("K3:K12").NumberFormat = "0.00"

# --- L column: Phi bei 1uF --------------------------------------------------
$ws.Range("L2").Formula = "=ATAN(2*PI()*50*G2*0.000001) * 180 / PI()"
$ws.Range("L2").NumberFormat = "0.00"
$ws.Range("L3:L12").Formula = "=ATAN(2*PI()*50*G3*0.000001) * 180 / PI()"
$ws.Range("L3:L12").NumberFormat = "0.00"

# --- M column: Ind fuer 20 Grad ---------------------------------------------
$ws.Range("M2").Formula = "= TAN(20 / 180 *PI()) * G2 / (2*PI()*50)"
$ws.Range("M2").NumberFormat = "0.0000"
$ws.Range("M3:M12").Formula = "= TAN(20 / 180 *PI()) * G3 /(2*PI()*50)"
$ws.Range("M3:M12").NumberFormat = "0.0000"

# --- N column: Phi bei 20mH -------------------------------------------------
$ws.Range("N2").Formula = "=ATAN(2*PI()*50*0.02/G2)"
$ws.Range("N2").NumberFormat = "0.000"
$ws.Range("N3:N12").Formula = "=ATAN(2*PI()*50*0.02/G3)"
$ws.Range("N3:N12").NumberFormat = "0.000"

# --- O column: Phi bei 10mH -------------------------------------------------
$ws.Range("O2").Formula = "=ATAN(2*PI()*50*0.01/G2)"
$ws.Range("O2").NumberFormat = "0.000"
$ws.Range("O3:O12").Formula = "=ATAN(2*PI()*50*0.01/G3)"
$ws.Range("O3:O12").NumberFormat = "0.000"

# --- Header row 13 -----------------------------------------------------------
$ws.Range("J13").Value = "Kap für (20°)"
$ws.Range("K13").Value = "Phi bei 10µF"
$ws.Range("L13").Value = "Phi bei 1µF"
$ws.Range("M13").Value = "Ind für 20°"
$ws.Range("N13").Value = "Phi bei 20mH"
$ws.Range("O13").Value = "Phi bei 10mH"

# --- Stray formatted (but empty) cell ----------------------------------------
$ws.Range("K14").NumberFormat = "0.00E+00"

# --- Column widths for the new columns ---------------------------------------
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666
$ws.Columns.Item(11).ColumnWidth = 12.833333333333334
$ws.Columns.Item(12).ColumnWidth = 11
$ws.Columns.Item(13).ColumnWidth = 10.666666666666666
$ws.Columns.Item(14).ColumnWidth = 12.833333333333334
$ws.Columns.Item(15).ColumnWidth = 12.166666666666666

# --- Move selection/scroll to K2 (drops the old topLeftCell="A28") -----------
$ws.Range("K2").Select() | Out-Null
